# Field Update : Code
# The header cell A1 ("ModuleName") is renamed to "Code".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "Code"
